$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 327.875
$ws.Range("I33").Value = 303.6154
$ws.Range("K33").Value = 303.6154
$ws.Range("M33").Value = -74.61540000000002

$ws.Range("H62").Value = 1919.6666
$ws.Range("I62").Value = 1164.2
$ws.Range("K62").Value = 1164.2
$ws.Range("M62").Value = -540.2

$ws.Range("H65").Value = 1919.6666
$ws.Range("I65").Value = 1164.2
$ws.Range("K65").Value = 5821
$ws.Range("M65").Value = -2701

$ws.Range("H86").Value = 1943.8572
$ws.Range("I86").Value = 2001
$ws.Range("J86").Value = 1901
$ws.Range("K86").Value = 2001
$ws.Range("L86").Value = 1901
$ws.Range("M86").Value = -878
$ws.Range("N86").Value = -4147

$ws.Range("H89").Value = 1943.8572
$ws.Range("I89").Value = 2001
$ws.Range("J89").Value = 1901
$ws.Range("K89").Value = 10005
$ws.Range("L89").Value = 9505
$ws.Range("M89").Value = -4389
$ws.Range("N89").Value = -20737

$ws.Range("H93").Value = 38900
$ws.Range("J93").Value = 38900
$ws.Range("L93").Value = 38900
$ws.Range("N93").Value = -43892

$ws.Range("H98").Value = 4868.9
$ws.Range("I98").Value = 1955.5714
$ws.Range("J98").Value = 11666.667
$ws.Range("K98").Value = 1955.5714
$ws.Range("L98").Value = 11666.667
$ws.Range("M98").Value = -457.5714
$ws.Range("N98").Value = -14662.667

$ws.Range("H122").Value = 4868.9
$ws.Range("I122").Value = 1955.5714
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 5866.7142
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -3416.7142
$ws.Range("N122").Value = -39900.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 67466.664
$ws.Range("J69").Value = 67466.664
$ws.Range("L69").Value = 67466.664
$ws.Range("N69").Value = -68964.664

$ws.Range("H72").Value = 67466.664
$ws.Range("J72").Value = 67466.664
$ws.Range("L72").Value = 202399.992
$ws.Range("N72").Value = -209887.992

$ws.Range("H97").Value = 1065.8928
$ws.Range("I97").Value = 711.7826
$ws.Range("K97").Value = 711.7826
$ws.Range("M97").Value = -215.7826

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1093.2727
$ws.Range("I94").Value = 1112.8889
$ws.Range("J94").Value = 1005
$ws.Range("K94").Value = 1112.8889
$ws.Range("L94").Value = 1005
$ws.Range("M94").Value = -661.8888999999999
$ws.Range("N94").Value = -1907

$ws.Range("H134").Value = 1950.907
$ws.Range("J134").Value = 7250
$ws.Range("L134").Value = 21750
$ws.Range("N134").Value = -26820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12502418
$ws.Range("I31").Value = 1213.6
$ws.Range("J31").Value = 33337760
$ws.Range("K31").Value = 1213.6
$ws.Range("L31").Value = 33337760
$ws.Range("M31").Value = -918.5999999999999
$ws.Range("N31").Value = -33338350

$ws.Range("H34").Value = 12502418
$ws.Range("I34").Value = 1213.6
$ws.Range("J34").Value = 33337760
$ws.Range("K34").Value = 1213.6
$ws.Range("L34").Value = 33337760
$ws.Range("M34").Value = -1011.6
$ws.Range("N34").Value = -33338164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 84.30769
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 88.44444
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 530.66664
$ws.Range("M2").Value = -337
$ws.Range("N2").Value = -756.66664

$ws.Range("H22").Value = 2418.3635
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2418.3635
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 7255.0905
$ws.Range("N22").Value = -7593.0905
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 2418.3635
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2418.3635
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 7255.0905
$ws.Range("N27").Value = -7459.0905
$ws.Range("M27").ClearContents()

$ws.Range("H47").Value = 990
$ws.Range("I47").Value = 990
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 2970
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = -2539
$ws.Range("M47").ClearContents()

$ws.Range("H50").Value = 413.83334
$ws.Range("I50").Value = 292.25
$ws.Range("J50").Value = 657
$ws.Range("K50").Value = 876.75
$ws.Range("L50").Value = 1971
$ws.Range("M50").Value = -395.75
$ws.Range("N50").Value = -2933

$ws.Range("H53").Value = 413.83334
$ws.Range("I53").Value = 292.25
$ws.Range("J53").Value = 657
$ws.Range("K53").Value = 876.75
$ws.Range("L53").Value = 1971
$ws.Range("M53").Value = -395.75
$ws.Range("N53").Value = -2933

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3225.9387
$ws.Range("I126").Value = 2899.2927
$ws.Range("J126").Value = 4900
$ws.Range("K126").Value = 8697.8781
$ws.Range("L126").Value = 14700
$ws.Range("M126").Value = -6227.8781
$ws.Range("N126").Value = -19640

$ws.Range("H132").Value = 2734.5881
$ws.Range("I132").Value = 1669.15
$ws.Range("K132").Value = 5007.450000000001
$ws.Range("M132").Value = -2477.450000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5038.609
$ws.Range("I40").Value = 4332.7617
$ws.Range("J40").Value = 12450
$ws.Range("K40").Value = 4332.7617
$ws.Range("L40").Value = 12450
$ws.Range("M40").Value = -4196.7617
$ws.Range("N40").Value = -12722

$ws.Range("H81").Value = 56763.168
$ws.Range("J81").Value = 56763.168
$ws.Range("L81").Value = 56763.168
$ws.Range("N81").Value = -58759.168

$ws.Range("H84").Value = 56763.168
$ws.Range("J84").Value = 56763.168
$ws.Range("L84").Value = 170289.504
$ws.Range("N84").Value = -180273.504

$ws.Range("H87").Value = 33333.332
$ws.Range("J87").Value = 33333.332
$ws.Range("L87").Value = 33333.332
$ws.Range("N87").Value = -35579.332

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H90").Value = 33333.332
$ws.Range("J90").Value = 33333.332
$ws.Range("L90").Value = 99999.99600000001
$ws.Range("N90").Value = -111231.996

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H93").Value = 3587019.5
$ws.Range("I93").Value = 5557786
$ws.Range("J93").Value = 3808.0908
$ws.Range("K93").Value = 5557786
$ws.Range("L93").Value = 3808.0908
$ws.Range("M93").Value = -5556538
$ws.Range("N93").Value = -6304.0908

$ws.Range("H136").Value = 3871.5833
$ws.Range("I136").Value = 1554.5333
$ws.Range("J136").Value = 7733.3335
$ws.Range("K136").Value = 4663.5999
$ws.Range("L136").Value = 23200.0005
$ws.Range("M136").Value = -2113.5999
$ws.Range("N136").Value = -28300.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 219.14285
$ws.Range("I107").Value = 192.11765
$ws.Range("J107").Value = 334
$ws.Range("K107").Value = 576.35295
$ws.Range("L107").Value = 1002
$ws.Range("M107").Value = 1343.64705
$ws.Range("N107").Value = -4842

$ws.Range("H136").Value = 18568.625
$ws.Range("I136").Value = 19758.334
$ws.Range("J136").Value = 14999.5
$ws.Range("K136").Value = 59275.00199999999
$ws.Range("L136").Value = 44998.5
$ws.Range("M136").Value = -56725.00199999999
$ws.Range("N136").Value = -50098.5
